# MitsosBarton2006Ex314 - Strong Stationary generator (alpha_zero)
# Update numeric/expression values across several sheets.
#
# Many of the target values look like plain numbers ("1.1", "0.3", ...).
# If assigned directly, Excel auto-converts them to numeric cells, but the
# workbook stores them as text (shared-string) cells. To preserve that,
# we force a text number-format before the assignment and then clear the
# format again so no stray style is left behind on the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A2").Value = "-2.1 + x"
Set-TextValue $ws.Range("B2") "1.1"
Set-TextValue $ws.Range("D2") "0.74"
$ws.Range("A3").Value = "2.1 - x"
Set-TextValue $ws.Range("B3") "-3.1"
Set-TextValue $ws.Range("D3") "0.27"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-3.3000000000000003 + y"
Set-TextValue $ws.Range("B2") "2.3000000000000003"
Set-TextValue $ws.Range("D2") "0.22"
Set-TextValue $ws.Range("E2") "7.1"
Set-TextValue $ws.Range("F2") "0.3"
$ws.Range("A3").Value = "3.3000000000000007 - y"
Set-TextValue $ws.Range("B3") "-4.300000000000001"
Set-TextValue $ws.Range("D3") "0.66"
Set-TextValue $ws.Range("E3") "8.8"
Set-TextValue $ws.Range("F3") "3.1"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "2.1"
Set-TextValue $ws.Range("B2") "3.3000000000000003"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws.Range("A2") "-8.350000000000003"

# --- Vector_BF ---
# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"); Worksheets.Item(<name>) resolves
# case-insensitively and would hit the wrong one ("Vector_bf"), so this
# sheet must be addressed by its 1-based tab position (6) instead.
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-4.17"
Set-TextValue $ws.Range("A3") "-4.8999999999999995"
